$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 3
$ws.Range("A3").Value = "'1756651754825"
$ws.Range("A3").Style = "Normal"
$ws.Range("B3").Value = "'1756651364942"
$ws.Range("B3").Style = "Normal"
$ws.Range("C3").Value = "'Sudip Maharjan"
$ws.Range("C3").Style = "Normal"
$ws.Range("D3").Value = "'1755340190541"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "'Ravibhawan Ghar"
$ws.Range("E3").Style = "Normal"
$ws.Range("F3").Value = 0
$ws.Range("G3").Value = 0
$ws.Range("H3").Value = 50000
$ws.Range("I3").Value = 0
$ws.Range("J3").Value = 0
$ws.Range("K3").Value = 50000
$ws.Range("L3").Value = "'2025-08-24"
$ws.Range("L3").Style = "Normal"
$ws.Range("M3").Value = "'Bank Transfer"
$ws.Range("M3").Style = "Normal"
$ws.Range("N3").Value = "'"
$ws.Range("N3").Style = "Normal"
$ws.Range("O3").Value = "'"
$ws.Range("O3").Style = "Normal"
$ws.Range("P3").Value = "'Nabil"
$ws.Range("P3").Style = "Normal"
$ws.Range("Q3").Value = "'"
$ws.Range("Q3").Style = "Normal"
$ws.Range("R3").Value = "'security_deposit"
$ws.Range("R3").Style = "Normal"
$ws.Range("S3").Value = "'2025-08-31T14:49:14.825Z"
$ws.Range("S3").Style = "Normal"

# Row 4
$ws.Range("A4").Value = "'1756655343517"
$ws.Range("A4").Style = "Normal"
$ws.Range("B4").Value = "'1756651364942"
$ws.Range("B4").Style = "Normal"
$ws.Range("C4").Value = "'Sudip Maharjan"
$ws.Range("C4").Style = "Normal"
$ws.Range("D4").Value = "'1755340190541"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "'Ravibhawan Ghar"
$ws.Range("E4").Style = "Normal"
$ws.Range("F4").Value = "'2"
$ws.Range("F4").Style = "Normal"
$ws.Range("G4").Value = "'2025"
$ws.Range("G4").Style = "Normal"
$ws.Range("H4").Value = 50000
$ws.Range("I4").Value = 0
$ws.Range("J4").Value = 0
$ws.Range("K4").Value = 50000
$ws.Range("L4").Value = "'2025-08-31"
$ws.Range("L4").Style = "Normal"
$ws.Range("M4").Value = "'Cash"
$ws.Range("M4").Style = "Normal"
$ws.Range("N4").Value = "'"
$ws.Range("N4").Style = "Normal"
$ws.Range("O4").Value = "'"
$ws.Range("O4").Style = "Normal"
$ws.Range("P4").Value = "'"
$ws.Range("P4").Style = "Normal"
$ws.Range("Q4").Value = "'"
$ws.Range("Q4").Style = "Normal"
$ws.Range("R4").Value = "'rent"
$ws.Range("R4").Style = "Normal"
$ws.Range("S4").Value = "'2025-08-31T15:49:03.518Z"
$ws.Range("S4").Style = "Normal"

# Row 5
$ws.Range("A5").Value = "'1756655447320"
$ws.Range("A5").Style = "Normal"
$ws.Range("B5").Value = "'1756651364942"
$ws.Range("B5").Style = "Normal"
$ws.Range("C5").Value = "'Sudip Maharjan"
$ws.Range("C5").Style = "Normal"
$ws.Range("D5").Value = "'1755340190541"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "'Ravibhawan Ghar"
$ws.Range("E5").Style = "Normal"
$ws.Range("F5").Value = "'11"
$ws.Range("F5").Style = "Normal"
$ws.Range("G5").Value = "'2025"
$ws.Range("G5").Style = "Normal"
$ws.Range("H5").Value = 50000
$ws.Range("I5").Value = 0
$ws.Range("J5").Value = 0
$ws.Range("K5").Value = 50000
$ws.Range("L5").Value = "'2026-01-15"
$ws.Range("L5").Style = "Normal"
$ws.Range("M5").Value = "'Cash"
$ws.Range("M5").Style = "Normal"
$ws.Range("N5").Value = "'"
$ws.Range("N5").Style = "Normal"
$ws.Range("O5").Value = "'"
$ws.Range("O5").Style = "Normal"
$ws.Range("P5").Value = "'"
$ws.Range("P5").Style = "Normal"
$ws.Range("Q5").Value = "'"
$ws.Range("Q5").Style = "Normal"
$ws.Range("R5").Value = "'rent"
$ws.Range("R5").Style = "Normal"
$ws.Range("S5").Value = "'2025-08-31T15:50:47.320Z"
$ws.Range("S5").Style = "Normal"

